$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cases")

# Rename three test case identifiers in column A
$ws.Range("A5").Value = "queryNonExistentOrder"
$ws.Range("A11").Value = "postOrderWithNegativeQuantity"
$ws.Range("A12").Value = "postOrderWithInvalidId"

# Update the active selection to D18
$ws.Activate()
$ws.Range("D18").Select()
